$p = $ppt.ActivePresentation

# --- Slide 23: "Classification Tree - Credit Policy" ---
$s23 = $p.Slides.Item(23)
$shp23 = $s23.Shapes.Item(5)

# Resize/reposition the stats text box to make room for two extra bullet lines
$shp23.Top = 151.086617
$shp23.Height = 375.632835

# Insert "Precision" / "Recall" bullets after the "Specificity: 0.794" line
$tr23 = $shp23.TextFrame.TextRange
$cr = [char]13
$specificity23 = $tr23.Paragraphs(7, 1)
[void]$specificity23.InsertAfter($cr + "Precision: 0.952" + $cr + "Recall: 0.999")

# --- Slide 24: "Classification Tree - Not Fully Paid" ---
$s24 = $p.Slides.Item(24)
$shp24 = $s24.Shapes.Item(5)

# Resize/reposition the stats text box to make room for two extra bullet lines
$shp24.Top = 136.262993
$shp24.Height = 398.655442

# Insert "Precision" / "Recall" bullets after the "Specificity: 0.013" line
$tr24 = $shp24.TextFrame.TextRange
$specificity24 = $tr24.Paragraphs(8, 1)
[void]$specificity24.InsertAfter($cr + "Precision: 0.645" + $cr + "Recall: 0.013")
